# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (positioned between the existing
# "2021-Q4" and "总计" sheets) populated with fund-holding detail rows,
# and updates the "总计" (totals) sheet with a new summary row for the
# 2022-Q1 quarter (prepended above the existing 2021-Q4 summary row).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet right after "2021-Q4" and before "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add($null, $q4)
$q1.Name = "2022-Q1"

# Reuse the header row (styles + text) and the index column's style from
# the "2021-Q4" sheet so the new sheet matches the established layout.
$q4.Range("A1:H1").Copy($q1.Range("A1:H1"))
$q4.Range("A2:A4").Copy($q1.Range("A2:A4"))

# Columns B-G hold text values (fund code/name/scale/position/etc, some
# of which look numeric) - force a text format so they're stored as
# strings rather than being auto-coerced to numbers.
$q1.Range("B2:G4").NumberFormat = "@"

$q1.Cells.Item(2, 2).Value = "004616"
$q1.Cells.Item(2, 3).Value = "中欧电子信息产业沪港深股票A"
$q1.Cells.Item(2, 4).Value = "14.54"
$q1.Cells.Item(2, 5).Value = "92.26"
$q1.Cells.Item(2, 6).Value = "3.94"
$q1.Cells.Item(2, 7).Value = "0.5729"
$q1.Cells.Item(2, 8).Value = 7

$q1.Cells.Item(3, 2).Value = "005763"
$q1.Cells.Item(3, 3).Value = "中欧电子信息产业沪港深股票C"
$q1.Cells.Item(3, 4).Value = "7.73"
$q1.Cells.Item(3, 5).Value = "92.26"
$q1.Cells.Item(3, 6).Value = "3.94"
$q1.Cells.Item(3, 7).Value = "0.3046"
$q1.Cells.Item(3, 8).Value = 7

$q1.Cells.Item(4, 2).Value = "006157"
$q1.Cells.Item(4, 3).Value = "财通量化核心优选混合"
$q1.Cells.Item(4, 4).Value = "0.09"
$q1.Cells.Item(4, 5).Value = "92.85"
$q1.Cells.Item(4, 6).Value = "1.61"
$q1.Cells.Item(4, 7).Value = "0.0014"
$q1.Cells.Item(4, 8).Value = 4

# ---------------------------------------------------------------------
# 2. Update the "总计" sheet: insert a new summary row for 2022-Q1 above
#    the existing 2021-Q4 row, and bump that row's index value.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
# The insert can bleed formatting from neighbouring rows into the blank
# row it creates - start clean, then re-apply only the index column's
# style (copied from the row below, matching the sheet's existing look).
$total.Range("A2:D2").ClearFormats()
$total.Range("A3").Copy($total.Range("A2"))

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 0.88

# The old row (now shifted down to row 3) keeps its data but its index
# column needs to advance from 0 to 1.
$total.Cells.Item(3, 1).Value = 1
